# Spare_Parts_Inventory.xlsx update - 2025-08-29 09:42:06
# Adds 4 new MLT_TESTER spare-part rows (Purlin, Clamps, M8-Bolts, Nuts)
# right after "Three color light indicator" (row 21) and before
# "Compressed air unit", pushing every following row down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 22 (rows below shift down, formatting
# of the row above is carried onto the new rows, matching the existing
# MLT_TESTER block formatting).
$ws.Range("A22:A25").EntireRow.Insert()

# Populate the newly inserted rows.
$ws.Cells.Item(22, 1).Value = "MLT_TESTER"
$ws.Cells.Item(22, 2).Value = "Purlin"

$ws.Cells.Item(23, 1).Value = "MLT_TESTER"
$ws.Cells.Item(23, 2).Value = "Clamps"

$ws.Cells.Item(24, 1).Value = "MLT_TESTER"
$ws.Cells.Item(24, 2).Value = "M8-Bolts"

$ws.Cells.Item(25, 1).Value = "MLT_TESTER"
$ws.Cells.Item(25, 2).Value = "Nuts"

# Match the author's final selection recorded in the workbook.
$ws.Range("C27").Select()
